$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 1212
$ws.Range("K3").Value = 1135
$ws.Range("I4").Value = 1785
$ws.Range("K4").Value = 244
$ws.Range("K5").Value = 69
$ws.Range("K6").Value = 1478
$ws.Range("I7").Value = 26238
$ws.Range("K7").Value = 4138

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K2").Value = 15
$ws.Range("K3").Value = 12
$ws.Range("K7").Value = 55

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 71
$ws.Range("K6").Value = 82
$ws.Range("K7").Value = 243

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K6").Value = 11
$ws.Range("K7").Value = 77

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K2").Value = 48
$ws.Range("K3").Value = 61
$ws.Range("K4").Value = 8
$ws.Range("K6").Value = 47
$ws.Range("K7").Value = 166

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K2").Value = 23
$ws.Range("K3").Value = 20
$ws.Range("K7").Value = 70

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K3").Value = 42
$ws.Range("K6").Value = 50
$ws.Range("K7").Value = 135

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K3").Value = 26
$ws.Range("K7").Value = 108

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("K3").Value = 27
$ws.Range("K7").Value = 75

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K4").Value = 16
$ws.Range("K7").Value = 125
$ws.Range("K8").Value = 243
$ws.Range("K14").Value = 24
$ws.Range("K15").Value = 30
$ws.Range("K18").Value = 32
$ws.Range("K19").Value = 111
$ws.Range("K20").Value = 101
$ws.Range("K24").Value = 15
$ws.Range("K25").Value = 19
$ws.Range("K27").Value = 50
$ws.Range("K29").Value = 191
$ws.Range("K31").Value = 46
$ws.Range("K33").Value = 166
$ws.Range("K36").Value = 46
$ws.Range("K37").Value = 135
$ws.Range("K42").Value = 136
$ws.Range("K43").Value = 39
$ws.Range("K46").Value = 8
$ws.Range("K48").Value = 39
$ws.Range("K52").Value = 114
$ws.Range("K53").Value = 55
$ws.Range("K54").Value = 75
$ws.Range("I55").Value = 314
$ws.Range("K55").Value = 43
$ws.Range("K57").Value = 8
$ws.Range("K64").Value = 25
$ws.Range("K65").Value = 108
$ws.Range("K67").Value = 171
$ws.Range("K68").Value = 12
$ws.Range("K71").Value = 12
$ws.Range("K72").Value = 16
$ws.Range("K76").Value = 55
$ws.Range("K77").Value = 31
$ws.Range("K78").Value = 60
$ws.Range("K83").Value = 77
$ws.Range("K84").Value = 31
$ws.Range("K88").Value = 54
$ws.Range("K89").Value = 62
$ws.Range("K91").Value = 47
$ws.Range("K93").Value = 16
$ws.Range("K94").Value = 49
$ws.Range("K95").Value = 70
$ws.Range("K98").Value = 26
$ws.Range("K99").Value = 75
$ws.Range("I101").Value = 26238
$ws.Range("K101").Value = 4138

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K6").Value = 18
$ws.Range("K7").Value = 46

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K2").Value = 51
$ws.Range("K3").Value = 50
$ws.Range("K6").Value = 59
$ws.Range("K7").Value = 171

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("K3").Value = 10
$ws.Range("K7").Value = 31

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K3").Value = 24
$ws.Range("K7").Value = 75

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 52
$ws.Range("K3").Value = 57
$ws.Range("K7").Value = 191

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K2").Value = 12
$ws.Range("K7").Value = 39

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K2").Value = 36
$ws.Range("K3").Value = 36
$ws.Range("K6").Value = 30
$ws.Range("K7").Value = 111

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K6").Value = 29
$ws.Range("K7").Value = 55

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("K2").Value = 10
$ws.Range("K4").Value = 2
$ws.Range("K7").Value = 24

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K3").Value = 42
$ws.Range("K6").Value = 56
$ws.Range("K7").Value = 136

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K2").Value = 21
$ws.Range("K7").Value = 60

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("K2").Value = 18
$ws.Range("I4").Value = 13
$ws.Range("K6").Value = 15
$ws.Range("I7").Value = 314
$ws.Range("K7").Value = 43

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("K3").Value = 6
$ws.Range("K4").Value = 1
$ws.Range("K7").Value = 15

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("K2").Value = 4
$ws.Range("K7").Value = 8

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K2").Value = 18
$ws.Range("K3").Value = 17
$ws.Range("K7").Value = 47

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("K2").Value = 7
$ws.Range("K7").Value = 25

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K2").Value = 29
$ws.Range("K7").Value = 101

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K2").Value = 12
$ws.Range("K7").Value = 32

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K2").Value = 18
$ws.Range("K4").Value = 3
$ws.Range("K7").Value = 46

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("K3").Value = 4
$ws.Range("K7").Value = 16

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K3").Value = 37
$ws.Range("K7").Value = 125

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("K2").Value = 17
$ws.Range("K4").Value = 5
$ws.Range("K7").Value = 49

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("K3").Value = 8
$ws.Range("K7").Value = 19

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("K3").Value = 3
$ws.Range("K7").Value = 30

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("K6").Value = 19
$ws.Range("K7").Value = 26

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K6").Value = 30
$ws.Range("K7").Value = 54

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K2").Value = 7
$ws.Range("K7").Value = 62

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("K2").Value = 17
$ws.Range("K3").Value = 6
$ws.Range("K6").Value = 21
$ws.Range("K7").Value = 50

$ws = $wb.Worksheets.Item('North Park')
$ws.Range("K3").Value = 4
$ws.Range("K7").Value = 12

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("K2").Value = 2
$ws.Range("K7").Value = 8

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("K3").Value = 13
$ws.Range("K4").Value = 3
$ws.Range("K6").Value = 18
$ws.Range("K7").Value = 39

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("K6").Value = 2
$ws.Range("K7").Value = 12

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("K4").Value = 3
$ws.Range("K7").Value = 16

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("K2").Value = 16
$ws.Range("K7").Value = 31

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K2").Value = 27
$ws.Range("K3").Value = 23
$ws.Range("K5").Value = 4
$ws.Range("K6").Value = 56
$ws.Range("K7").Value = 114

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("K2").Value = 4
$ws.Range("K7").Value = 16
